$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row values (B2:E2)
$ws.Range("B2").Value = 34.377179995401043
$ws.Range("C2").Value = 43.514205037357755
$ws.Range("D2").Value = 38.247644724786852
$ws.Range("E2").Value = 41.788895704827226

# Update STR row values (B3:E3)
$ws.Range("B3").Value = 42.45096915661842
$ws.Range("C3").Value = 46.901860647813173
$ws.Range("D3").Value = 44.281495040670407
$ws.Range("E3").Value = 40.5078180312737

# Update selection to reflect new selected range B1:E3
$ws.Range("B1:E3").Select()
